# Update "camas_ocupadas" dataset: append two new daily rows (2020-08-31 and 2020-09-01)
# that were missing in the previous export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data to append, in order
$newRows = @(
    @{ Row = 113; Fecha = "2020-08-31"; Territorio = "Andalucía"; CamasOcup = 493; UCI = 72 },
    @{ Row = 114; Fecha = "2020-09-01"; Territorio = "Andalucía"; CamasOcup = 522; UCI = 70 }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    # Column A: date stored as text (matches the rest of the column, which uses a "@" text format)
    $cellA = $ws.Range("A$rowNum")
    $cellA.NumberFormat = "@"
    $cellA.Value = $r.Fecha

    # Column B: territory name (plain text)
    $cellB = $ws.Range("B$rowNum")
    $cellB.Value = $r.Territorio

    # Column C: occupied beds count (number)
    $cellC = $ws.Range("C$rowNum")
    $cellC.Value = $r.CamasOcup

    # Column D: ICU count (number)
    $cellD = $ws.Range("D$rowNum")
    $cellD.Value = $r.UCI
}

# Select the last entered cell, mirroring the author's final cursor position
[void]$ws.Range("D114").Select()
